$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 150
$ws.Range("I12").Value = 150
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 20

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5460
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 5000
$ws.Range("N35").Value = -5758
$ws.Range("M35").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5240.778
$ws.Range("I40").Value = 3239.4
$ws.Range("J40").Value = 7742.5
$ws.Range("K40").Value = 3239.4
$ws.Range("L40").Value = 7742.5
$ws.Range("M40").Value = -3064.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3128
$ws.Range("I98").Value = 721.34375
$ws.Range("J98").Value = 11685
$ws.Range("K98").Value = 721.34375
$ws.Range("L98").Value = 11685
$ws.Range("M98").Value = 776.65625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3128
$ws.Range("I122").Value = 721.34375
$ws.Range("J122").Value = 11685
$ws.Range("K122").Value = 2164.03125
$ws.Range("L122").Value = 35055
$ws.Range("M122").Value = 285.96875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1519.52
$ws.Range("I132").Value = 1528.9778
$ws.Range("J132").Value = 1434.4
$ws.Range("K132").Value = 4586.9334
$ws.Range("L132").Value = 4303.200000000001
$ws.Range("M132").Value = -2056.9334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 31253358
$ws.Range("I137").Value = 71431070
$ws.Range("J137").Value = 4024
$ws.Range("K137").Value = 214293210
$ws.Range("L137").Value = 12072
$ws.Range("M137").Value = -214290660
$ws.Range("N137").Value = -17172

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2966.3333
$ws.Range("I138").Value = 2598.923
$ws.Range("J138").Value = 3131.0344
$ws.Range("K138").Value = 7796.768999999999
$ws.Range("L138").Value = 9393.1032
$ws.Range("M138").Value = -2656.768999999999
$ws.Range("N138").Value = -19673.1032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24393898
$ws.Range("I32").Value = 25644048
$ws.Range("J32").Value = 16007
$ws.Range("K32").Value = 25644048
$ws.Range("L32").Value = 16007
$ws.Range("M32").Value = -25643761

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1590.6666
$ws.Range("I97").Value = 1758.875
$ws.Range("J97").Value = 245
$ws.Range("K97").Value = 1758.875
$ws.Range("L97").Value = 245
$ws.Range("M97").Value = -1262.875
$ws.Range("N97").Value = -1237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 7866.909
$ws.Range("I110").Value = 5949.5
$ws.Range("J110").Value = 8962.571
$ws.Range("K110").Value = 5949.5
$ws.Range("L110").Value = 8962.571
$ws.Range("M110").Value = -3904.5
$ws.Range("N110").Value = -13052.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 5008.5
$ws.Range("I24").Value = 3350
$ws.Range("J24").Value = 9984
$ws.Range("K24").Value = 3350
$ws.Range("L24").Value = 9984
$ws.Range("M24").Value = -3115

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11794.259
$ws.Range("I105").Value = 12477.368
$ws.Range("J105").Value = 10171.875
$ws.Range("K105").Value = 12477.368
$ws.Range("L105").Value = 10171.875
$ws.Range("M105").Value = -10730.368
$ws.Range("N105").Value = -13665.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 673.5
$ws.Range("I7").Value = 930.9091
$ws.Range("J7").Value = 107.2
$ws.Range("K7").Value = 930.9091
$ws.Range("L7").Value = 107.2
$ws.Range("M7").Value = -817.9091
$ws.Range("N7").Value = -333.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5148.6665
$ws.Range("I16").Value = 1965.6666
$ws.Range("J16").Value = 8331.666999999999
$ws.Range("K16").Value = 1965.6666
$ws.Range("L16").Value = 8331.666999999999
$ws.Range("M16").Value = -1678.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7736.4287
$ws.Range("I62").Value = 4487.25
$ws.Range("J62").Value = 12068.667
$ws.Range("K62").Value = 4487.25
$ws.Range("L62").Value = 12068.667
$ws.Range("M62").Value = -3863.25
$ws.Range("N62").Value = -13316.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7736.4287
$ws.Range("I65").Value = 4487.25
$ws.Range("J65").Value = 12068.667
$ws.Range("K65").Value = 22436.25
$ws.Range("L65").Value = 60343.335
$ws.Range("M65").Value = -19316.25
$ws.Range("N65").Value = -66583.33499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 5148.6665
$ws.Range("I113").Value = 1965.6666
$ws.Range("J113").Value = 8331.666999999999
$ws.Range("K113").Value = 1965.6666
$ws.Range("L113").Value = 8331.666999999999
$ws.Range("M113").Value = 204.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 39629.83
$ws.Range("I122").Value = 58424.668
$ws.Range("J122").Value = 8874.637000000001
$ws.Range("K122").Value = 175274.004
$ws.Range("L122").Value = 26623.911
$ws.Range("M122").Value = -172824.004
$ws.Range("N122").Value = -31523.911

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7397.5713
$ws.Range("I132").Value = 7260
$ws.Range("J132").Value = 7500.75
$ws.Range("K132").Value = 21780
$ws.Range("L132").Value = 22502.25
$ws.Range("M132").Value = -19250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 66668770
$ws.Range("I107").Value = 2593.8333
$ws.Range("J107").Value = 111112890
$ws.Range("K107").Value = 7781.499899999999
$ws.Range("L107").Value = 333338670
$ws.Range("M107").Value = -5861.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -16064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10199.462
$ws.Range("I70").Value = 7627.467
$ws.Range("J70").Value = 13706.728
$ws.Range("K70").Value = 7627.467
$ws.Range("L70").Value = 13706.728
$ws.Range("M70").Value = -7357.467
$ws.Range("N70").Value = -14246.728

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10199.462
$ws.Range("I73").Value = 7627.467
$ws.Range("J73").Value = 13706.728
$ws.Range("K73").Value = 7627.467
$ws.Range("L73").Value = 13706.728
$ws.Range("M73").Value = -6691.467
$ws.Range("N73").Value = -15578.728

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20348
$ws.Range("I93").Value = 17228
$ws.Range("J93").Value = 21388
$ws.Range("K93").Value = 17228
$ws.Range("L93").Value = 21388
$ws.Range("M93").Value = -15356
$ws.Range("N93").Value = -25132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6230
$ws.Range("I68").Value = 4399.6665
$ws.Range("J68").Value = 7014.4287
$ws.Range("K68").Value = 4399.6665
$ws.Range("L68").Value = 7014.4287
$ws.Range("M68").Value = -3650.6665
$ws.Range("N68").Value = -8512.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 6230
$ws.Range("I71").Value = 4399.6665
$ws.Range("J71").Value = 7014.4287
$ws.Range("K71").Value = 21998.3325
$ws.Range("L71").Value = 35072.14350000001
$ws.Range("M71").Value = -18254.3325
$ws.Range("N71").Value = -42560.14350000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8488
$ws.Range("I122").Value = 4797.6665
$ws.Range("J122").Value = 10333.167
$ws.Range("K122").Value = 14392.9995
$ws.Range("L122").Value = 30999.501
$ws.Range("M122").Value = -11942.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 36000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 36000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 36000
$ws.Range("N47").Value = -37144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 50000
$ws.Range("I57").Value = 50000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 50000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -49246

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23796.8
$ws.Range("I64").Value = 14502.4
$ws.Range("J64").Value = 33091.2
$ws.Range("K64").Value = 14502.4
$ws.Range("L64").Value = 33091.2
$ws.Range("M64").Value = -14254.4
$ws.Range("N64").Value = -33587.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 23796.8
$ws.Range("I67").Value = 14502.4
$ws.Range("J67").Value = 33091.2
$ws.Range("K67").Value = 14502.4
$ws.Range("L67").Value = 33091.2
$ws.Range("M67").Value = -13644.4
$ws.Range("N67").Value = -34807.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 27124.5
$ws.Range("I70").Value = 28332.834
$ws.Range("J70").Value = 23499.5
$ws.Range("K70").Value = 28332.834
$ws.Range("L70").Value = 23499.5
$ws.Range("M70").Value = -28017.834
$ws.Range("N70").Value = -24129.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 27124.5
$ws.Range("I73").Value = 28332.834
$ws.Range("J73").Value = 23499.5
$ws.Range("K73").Value = 28332.834
$ws.Range("L73").Value = 23499.5
$ws.Range("M73").Value = -27240.834
$ws.Range("N73").Value = -25683.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 49997
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 49997
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 49997
$ws.Range("N93").Value = -54989

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 30000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 30000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 660.5
$ws.Range("I113").Value = 358.2
$ws.Range("J113").Value = 876.4286
$ws.Range("K113").Value = 1074.6
$ws.Range("L113").Value = 2629.2858
$ws.Range("M113").Value = 1095.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4170.1333
$ws.Range("I122").Value = 1887.8462
$ws.Range("J122").Value = 19005
$ws.Range("K122").Value = 5663.5386
$ws.Range("L122").Value = 57015
$ws.Range("M122").Value = -3213.5386
$ws.Range("N122").Value = -61915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5923.0527
$ws.Range("I132").Value = 2417.375
$ws.Range("J132").Value = 11932.786
$ws.Range("K132").Value = 7252.125
$ws.Range("L132").Value = 35798.358
$ws.Range("M132").Value = -4722.125
